$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A17").Value = "Judith"
$ws.Range("B17").Value = 650
